$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column (Price) and other text cells keep their literal text representation
# (avoids Excel auto-converting numeric-looking strings into floating point numbers).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.138.66'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.475.99'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.48%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '563.30'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '163.30'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.91%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.508'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.473.59'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.153'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.97%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.165'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.334'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.57%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.88'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '69.041.54'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.22%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000170'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '23.77'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.477.81'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.62%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.82'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.90%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '340.39'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.07'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.81'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.92'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +4.42%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '67.49'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.71'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.610.56'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.81%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.27'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0826'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.60%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.23'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '436.35'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.15'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.85%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.63'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.71%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '157.45'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.06'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.41%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.66%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.88'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.11%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'PolygonEcosystemToken'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.303'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.52%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.47'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '37.64'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.49'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.24%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.09'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '133.90'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.82%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.36'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0718'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.487'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.565'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.56%  '
